$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: tiny precision corrections to existing numeric values ---
$ws.Range("B2").Value = 0.04215534119371416
$ws.Range("D2").Value = 0.1361288253571671

# --- Row 1: add two new header cells (G1, H1), reusing the same
#     formatting (bold, centered, bordered) already used by the other
#     header cells, by copying the format from F1 ---
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# --- Row 2: add the new data values under the new headers ---
$ws.Range("G2").Value = 0.1194315095165318
$ws.Range("H2").Value = 0.9890000000000001
